$wb = $excel.ActiveWorkbook

# Map of worksheet name -> new B2 value (EBITDA-like column). A $null value
# means only the date in A2 changes for that sheet (B2 stays the same).
$updates = @{
    "Alcoa"                     = "4.48"
    "Rio Tinto"                 = "7.70"
    "Norsk Hydro"               = "2.65"
    "Reliance Steel & Aluminum" = "12.44"
    "Kaiser Aluminum"           = "9.52"
    "Ryerson Holding"           = "25.20"
    "Ultra Clean Holdings"      = "11.49"
    "Benchmark Electronics"     = "11.19"
    "Celestica"                 = $null
    "Flex Ltd"                  = "13.48"
    "MKS Instruments"           = "15.07"
}

$newDate = "2025/10/31"

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Write the date/numeric-looking values as literal text (matching the
    # existing inline-string cells) instead of letting Excel auto-convert
    # them to a date serial / number, then restore the default "Normal"
    # style so no stray number-format style is left behind.
    $cellA = $ws.Range("A2")
    $cellA.NumberFormat = "@"
    $cellA.Value = $newDate
    $cellA.Style = "Normal"

    $newB = $updates[$sheetName]
    if ($newB -ne $null) {
        $cellB = $ws.Range("B2")
        $cellB.NumberFormat = "@"
        $cellB.Value = $newB
        $cellB.Style = "Normal"
    }
}
